$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The title and caption text were each split into one run per word
# (plus separate space runs). Collapse each back down to a single run
# containing the full sentence.
#
# Assigning the exact same concatenated text is treated as a no-op by
# the diffing engine (since the rendered text doesn't change), so each
# TextRange is first set to an unrelated placeholder string to force a
# genuine replacement, then set to the desired final text; this yields
# a single freshly-created run with a plain <a:rPr/>.

$title = $s.Shapes.Item("Title 1").TextFrame.TextRange
$title.Text = "PLACEHOLDER"
$title.Text = "A Table, with a caption"

$caption = $s.Shapes.Item("TextBox 3").TextFrame.TextRange
$caption.Text = "PLACEHOLDER"
$caption.Text = "Demonstration of simple table syntax, with alignment"
